# daily auto push: 2026-03-01 05:02 UTC
# Insert one new data row at row 899 (pushing the existing rows 899-940 down
# to 900-941) and populate the new row with the newly logged entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 899..940 down to 900..941, leaving a blank row at 899.
$ws.Rows.Item(899).Insert()

# Fill the newly inserted row with the new log entry. The date column
# stores plain text (e.g. "2026/12/29") rather than a real date value
# throughout this sheet, so force a text format before writing it --
# otherwise Excel auto-converts the "2026/03/01"-looking string into a
# date serial number. ClearFormats() afterwards drops the temporary
# text-format style so the cell matches the unstyled look of its peers.
$ws.Range("A899").NumberFormat = "@"
$ws.Range("A899").Value = "2026/03/01"
$ws.Range("A899").ClearFormats()

$ws.Range("B899").Value = "日"
$ws.Range("C899").Value = 13
$ws.Range("D899").Value = 201
